# Update Name of Algo
# Applies updated numeric values to a handful of cells in Sheet1,
# matching the target commit's data changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = -6.675199999999989
$ws.Range("B3").Value = 5.850199999999988
$ws.Range("D5").Value = -8.519099999999991
$ws.Range("B14").Value = 8.975800000000003
$ws.Range("B21").Value = 5.640999999999996
$ws.Range("B23").Value = 5.816099999999999
$ws.Range("B25").Value = 5.870099999999993
